# Reporte Por Unidad: update per-unit call quality breakdown with new
# figures from Anny's report (unit 22 split into two rows, unit 23 gains
# an "Inválida" row, new subtotal/grand-total rows appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Propagate row-level formatting (styles) to their new destinations
#    BEFORE any values are overwritten. Order matters: rows that act as
#    both a copy source and a later destination are handled source-first.
# ---------------------------------------------------------------------

# Old row 12 ("Total:" grand total) moves its look to the new row 13
# FIRST, before row 12 itself is restyled in the next step.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial($xlPasteFormats)

# Old row 11 ("Total por unidad:" for unit 23) moves its look to the new
# row 12 (new "Total por unidad:" position for unit 23).
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial($xlPasteFormats)

# Detail-row look (row 4) is reused for the new rows 7, 8 and 11.
$ws.Range("A4:E4").Copy()
$ws.Range("A7:E7").PasteSpecial($xlPasteFormats)
$ws.Range("A4:E4").Copy()
$ws.Range("A8:E8").PasteSpecial($xlPasteFormats)
$ws.Range("A4:E4").Copy()
$ws.Range("A11:E11").PasteSpecial($xlPasteFormats)

# Subtotal-row look (row 6) is reused for the new row 9.
$ws.Range("A6:E6").Copy()
$ws.Range("A9:E9").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Clear stale leftover content so cells that must end up blank
#    (the "A"/"D" columns on subtotal/total rows) don't keep old values.
# ---------------------------------------------------------------------
$ws.Range("A6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("A9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("A12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("D13").ClearContents()

# ---------------------------------------------------------------------
# 3) Write the new values, row by row. Percent cells are written last
#    (after all formatting copies) so the forced-text number format
#    sticks instead of being clobbered by a later PasteSpecial.
# ---------------------------------------------------------------------

function Set-PercentText($cell, $text) {
    # Force literal text (not an auto-converted percentage number).
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 4: Unidad 21 / Inválida
$ws.Range("A4").Value = 21
$ws.Range("B4").Value = "Inválida"
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 2000

# Row 5: Unidad 21 / Válida
$ws.Range("A5").Value = 21
$ws.Range("B5").Value = "Válida"
$ws.Range("C5").Value = 2
$ws.Range("E5").Value = 3100

# Row 6: Unidad 21 / Total por unidad
$ws.Range("B6").Value = "Total por unidad:"
$ws.Range("C6").Value = 4
$ws.Range("E6").Value = 5100

# Row 7: Unidad 22 / Inválida (new row)
$ws.Range("A7").Value = 22
$ws.Range("B7").Value = "Inválida"
$ws.Range("C7").Value = 1
$ws.Range("E7").Value = 365

# Row 8: Unidad 22 / Válida (new row)
$ws.Range("A8").Value = 22
$ws.Range("B8").Value = "Válida"
$ws.Range("C8").Value = 2
$ws.Range("E8").Value = 850

# Row 9: Unidad 22 / Total por unidad
$ws.Range("B9").Value = "Total por unidad:"
$ws.Range("C9").Value = 3
$ws.Range("E9").Value = 1215

# Row 10: Unidad 23 / Válida
$ws.Range("A10").Value = 23
$ws.Range("B10").Value = "Válida"
$ws.Range("C10").Value = 5
$ws.Range("E10").Value = 4350

# Row 11: Unidad 23 / Inválida (new row)
$ws.Range("A11").Value = 23
$ws.Range("B11").Value = "Inválida"
$ws.Range("C11").Value = 5
$ws.Range("E11").Value = 4455

# Row 12: Unidad 23 / Total por unidad (new position)
$ws.Range("B12").Value = "Total por unidad:"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = 8805

# Row 13: Grand total (new position)
$ws.Range("B13").Value = "Total:"
$ws.Range("C13").Value = 17
$ws.Range("E13").Value = 15120

# Percent (text) cells of column D, set last so the forced text format
# is not overwritten by any later formatting operation.
Set-PercentText $ws.Range("D4") "50%"
Set-PercentText $ws.Range("D5") "50%"
Set-PercentText $ws.Range("D7") "33%"
Set-PercentText $ws.Range("D8") "67%"
Set-PercentText $ws.Range("D10") "50%"
Set-PercentText $ws.Range("D11") "50%"
